$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Configs")

# 1) Insert a brand-new row above the existing header row (row 3).
#    This pushes the header + all the year groups down by one row and
#    shifts the existing merged B-column cells (B4:B6 -> B5:B7, etc.)
#    automatically.
$ws.Rows("3:3").Insert()

# Row 3 inherits the grey "section header" look from the row it was
# pushed down from -- clear that back to a plain, unformatted look before
# styling it the way the new "properties" row should look.
$ws.Range("B3:E3").ClearFormats()

# 2) Populate the freshly inserted row with the new "properties" labels.
$ws.Range("B3").Value = "properties"
$ws.Range("C3").Value = "origin"
$ws.Range("D3").Value = "Deviation"

# Keep the look-and-feel close to the rest of the unformatted sheet
# (plain Calibri, no fill) and add a thin right border to D3 to set the
# new "Deviation" column apart, matching the accent border on the new row.
$ws.Range("B3:E3").Font.Name = "Calibri"
$ws.Range("B3:E3").Font.Size = 11
$ws.Range("D3").Borders.Item(10).LineStyle = 1
$ws.Range("D3").Borders.Item(10).Weight = 2

# 3) Append a new blank row after the last data row (old row 18, now row
#    19), giving the table a trailing spacer row.
$ws.Range("B20:E20").ClearFormats()
$ws.Range("B20:E20").Font.Name = "Calibri"
$ws.Range("B20:E20").Font.Size = 11
